# Update "paises.xlsx" - update country rankings & case numbers (provincias Spain update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 02:05"

# --- Row 4: Estados Unidos - refreshed case counts ---
$ws.Range("B4").Value = 1644879
$ws.Range("C4").Value = 23982
$ws.Range("D4").Value = 397653
$ws.Range("E4").Value = 1149591
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1281
$ws.Range("H4").Value = 97635

# --- Row 17: Canada - refreshed case counts ---
$ws.Range("B17").Value = 82480
$ws.Range("C17").Value = 1156
$ws.Range("D17").Value = 42594
$ws.Range("E17").Value = 33636
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 98
$ws.Range("H17").Value = 6250

# --- Rows 42-43: Japon overtakes Austria in ranking ---
$ws.Range("A42").Value = "Japon"
$ws.Range("B42").Value = 16513
$ws.Range("C42").Value = 89
$ws.Range("D42").Value = 13005
$ws.Range("E42").Value = 2712
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 19
$ws.Range("H42").Value = 796

$ws.Range("A43").Value = "Austria"
$ws.Range("B43").Value = 16436
$ws.Range("C43").Value = 32
$ws.Range("D43").Value = 15005
$ws.Range("E43").Value = 796
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 635

# --- Rows 112-116: Venezuela jumps ahead of Niger / Chipre / Zambia / Costa Rica ---
$ws.Range("A112").Value = "Venezuela"
$ws.Range("B112").Value = 944
$ws.Range("C112").Value = 62
$ws.Range("D112").Value = 262
$ws.Range("E112").Value = 672
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 10

$ws.Range("A113").Value = "Niger"
$ws.Range("B113").Value = 937
$ws.Range("C113").Value = 13
$ws.Range("D113").Value = 764
$ws.Range("E113").Value = 113
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 60

$ws.Range("A114").Value = "Republica de Chipre"
$ws.Range("B114").Value = 927
$ws.Range("C114").Value = 4
$ws.Range("D114").Value = 561
$ws.Range("E114").Value = 349
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 17

$ws.Range("A115").Value = "Zambia"
$ws.Range("B115").Value = 920
$ws.Range("C115").Value = 54
$ws.Range("D115").Value = 336
$ws.Range("E115").Value = 577
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 7

$ws.Range("A116").Value = "Costa Rica"
$ws.Range("B116").Value = 911
$ws.Range("C116").Value = 8
$ws.Range("D116").Value = 600
$ws.Range("E116").Value = 301
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 10

# --- Row 144: Vietnam - refreshed case counts ---
$ws.Range("D144").Value = 267
$ws.Range("E144").Value = 57

# --- Rows 148-149: Guayana Francesa overtakes Santo Tome y Principe ---
$ws.Range("A148").Value = "Guayana Francesa"
$ws.Range("B148").Value = 261
$ws.Range("C148").Value = 12
$ws.Range("D148").Value = 141
$ws.Range("E148").Value = 119
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 1

$ws.Range("A149").Value = "Santo Tome y Principe"
$ws.Range("B149").Value = 251
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 4
$ws.Range("E149").Value = 239
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 8

# --- Row 159: Guadalupe - refreshed case counts ---
$ws.Range("B159").Value = 156
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 115
$ws.Range("E159").Value = 28

# --- Rows 209-210: Seychelles / Groenlandia swap places (tied counts) ---
$ws.Range("A209").Value = "Seychelles"
$ws.Range("A210").Value = "Groenlandia"

# --- Rows 214-216: Sahara Occidental / San Bartolome move ahead of Bonaire, San Eustaquio y Saba ---
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
